$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 40-41; this pushes the existing rows
# 40-53 down to 42-55, preserving all of their data untouched.
$ws.Rows("40:41").Insert()

# Row 40: new weekly observation (Provincia de Limari, Primera)
$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C40").Value = "Ñuble"
$ws.Range("D40").Value = 44837
$ws.Range("E40").Value = 16
$ws.Range("F40").Value = 100112026
$ws.Range("G40").Value = "Haba"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 30
$ws.Range("K40").Value = 9000
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = 9000
$ws.Range("N40").Value = "$/saco 25 kilos"
$ws.Range("O40").Value = "Provincia de Limarí"
$ws.Range("P40").Value = 360
$ws.Range("Q40").Value = 25
$ws.Range("R40").Value = "Hortaliza"

# Row 41: new weekly observation (Provincia de Limari, Segunda)
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44837
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112026
$ws.Range("G41").Value = "Haba"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Segunda"
$ws.Range("J41").Value = 30
$ws.Range("K41").Value = 9500
$ws.Range("L41").Value = 9500
$ws.Range("M41").Value = 9500
$ws.Range("N41").Value = "$/saco 25 kilos"
$ws.Range("O41").Value = "Provincia de Limarí"
$ws.Range("P41").Value = 380
$ws.Range("Q41").Value = 25
$ws.Range("R41").Value = "Hortaliza"

# Apply the date-number-format style (style index 2 in the original
# workbook) to the Fecha cells of the new rows, matching the other
# rows in column D.
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D41").NumberFormat = "YYYY-MM-DD HH:MM:SS"
